# Update the cached "datetimeFigureOut" date field text that appears on
# every slide layout's and the slide master's date placeholder, plus fix
# the casing of the "Spagi_76_Chairs_v_2" file/URL references on slide 3
# ("v" -> "V").

$p = $ppt.ActivePresentation

$oldDate = "21.04.2025"
$newDate = "23.04.2025"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

# Every slide layout's date placeholder.
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $cl = $master.CustomLayouts.Item($j)
    Update-DatePlaceholders $cl.Shapes
}

# Fix the "Spagi_76_Chairs_v_2" / "Spagi_76_Chairs_v_2.xml" text runs on
# slide 3 (lower-case "v" -> upper-case "V"). Rewriting a character range
# that ends in the *middle* of a run splits that run in two, so the
# replacement range is always extended through to the end of the text
# (re-appending whatever already followed the match unchanged) to keep
# each run intact.
$oldName = "Spagi_76_Chairs_v_2"
$newName = "Spagi_76_Chairs_V_2"

$s3 = $p.Slides.Item(3)
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $sh = $s3.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        $full = $tr.Text
        $pos = $full.IndexOf($oldName)
        if ($pos -ge 0) {
            $start = $pos + 1
            $restLen = $tr.Length - $pos
            $tail = $full.Substring($pos + $oldName.Length)
            $sub = $tr.Characters($start, $restLen)
            $sub.Text = $newName + $tail
        }
    }
}
